# Refresh the cryptocurrency price / 1h-volume snapshot (GitHub Actions data pull).
# Rows 21/22 and rows 35/36 also swap position (ranking reorder) as part of the refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Price" (column D) cells are stored as literal text in the source file. Prefixing the
# new value with an apostrophe forces Excel to keep it as text too, instead of silently
# coercing it to a number and dropping formatting (e.g. "59.30" -> 59.3, "0.100" -> 0.1,
# "1.00" -> 1).
$apos = "'"

$ws.Range('D2').Value = $apos + '36.764.10'
$ws.Range('E2').Value = '  +0.90%  '
$ws.Range('D3').Value = $apos + '1.967.72'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = $apos + '244.78'
$ws.Range('E5').Value = '  +0.57%  '
$ws.Range('D6').Value = $apos + '0.621'
$ws.Range('E6').Value = '  +1.03%  '
$ws.Range('D7').Value = $apos + '59.30'
$ws.Range('E7').Value = '  +2.51%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +2.18%  '
$ws.Range('D10').Value = $apos + '0.0820'
$ws.Range('E10').Value = '  -2.31%  '
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').Value = $apos + '22.42'
$ws.Range('E12').Value = '  +3.69%  '
$ws.Range('D13').Value = $apos + '2.255.88'
$ws.Range('E13').Value = '  +1.26%  '
$ws.Range('D14').Value = $apos + '0.832'
$ws.Range('E14').Value = '  +0.93%  '
$ws.Range('D15').Value = $apos + '13.81'
$ws.Range('E15').Value = '  +1.50%  '
$ws.Range('D16').Value = $apos + '5.29'
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('D17').Value = $apos + '1.969.28'
$ws.Range('E17').Value = '  +1.25%  '
$ws.Range('D18').Value = $apos + '36.660.91'
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('D19').Value = $apos + '69.88'
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').Value = $apos + '0.0₃0863'
$ws.Range('E20').Value = '  -0.66%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = $apos + '5.10'
$ws.Range('E21').Value = '  +1.18%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value = $apos + '229.40'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('E25').Value = '  +3.05%  '
$ws.Range('D26').Value = $apos + '9.35'
$ws.Range('E26').Value = '  +0.87%  '
$ws.Range('E27').Value = '  +15.16%  '
$ws.Range('D28').Value = $apos + '160.82'
$ws.Range('E28').Value = '  -0.88%  '
$ws.Range('D29').Value = $apos + '19.44'
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').Value = $apos + '0.120'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('D31').Value = $apos + '1.14'
$ws.Range('E31').Value = '  -1.13%  '
$ws.Range('D32').Value = $apos + '4.73'
$ws.Range('E32').Value = '  +1.20%  '
$ws.Range('D33').Value = $apos + '0.0622'
$ws.Range('E33').Value = '  -1.19%  '
$ws.Range('E34').Value = '  +0.64%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = $apos + '2.28'
$ws.Range('E35').Value = '  +6.48%  '
$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').Value = $apos + '1.00'
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('E37').Value = '  -2.56%  '
$ws.Range('D38').Value = $apos + '3.40'
$ws.Range('E38').Value = '  +11.68%  '
$ws.Range('E39').Value = '  -0.13%  '
$ws.Range('D40').Value = $apos + '0.0999'
$ws.Range('E40').Value = '  +2.66%  '
$ws.Range('E41').Value = '  -1.90%  '
$ws.Range('E42').Value = '  +1.79%  '
$ws.Range('E43').Value = '  -0.56%  '
$ws.Range('D44').Value = $apos + '16.16'
$ws.Range('E44').Value = '  +0.60%  '
$ws.Range('D45').Value = $apos + '1.365.81'
$ws.Range('E45').Value = '  +0.81%  '
$ws.Range('E46').Value = '  +0.99%  '
$ws.Range('D47').Value = $apos + '88.07'
$ws.Range('E47').Value = '  +0.37%  '
$ws.Range('D48').Value = $apos + '7.18'
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('E49').Value = '  +0.93%  '
$ws.Range('D50').Value = $apos + '2.146.18'
$ws.Range('E50').Value = '  +1.19%  '
$ws.Range('D51').Value = $apos + '44.12'
$ws.Range('E51').Value = '  -2.31%  '
